$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D and E (shifts old D:I -> F:K)
$ws.Range("D1:E1").EntireColumn.Insert()

# New header cells
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# New data cells for existing row 2
$ws.Range("D2").Value = "Mumbai Indians"
$ws.Range("E2").Value = "Rajasthan Royals"

# New row 3 with a fresh record
$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " October 28 2020"
$ws.Range("C3").Value = "Mumbai won by 5 wickets (with 5 balls remaining)"
$ws.Range("D3").Value = "Mumbai Indians"
$ws.Range("E3").Value = "Royal Challengers Bangalore"
$ws.Range("F3").Value = "Kieron Pollard (c)"
$ws.Range("G3").Value = "4"
$ws.Range("H3").Value = "1"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "400.00"
